$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values scraped from coinranking.com for the GitHub Actions refresh.
# Values that look like plain decimal numbers (e.g. "216.43") must be forced
# to remain text (matching the original inlineStr cells) instead of being
# auto-converted to numbers by Excel's smart cell-entry parsing, and then the
# temporary text number-format is reverted so the cell keeps its original
# (default) style.

$ws.Range("D2").Value = "26.032.75"
$ws.Range("D3").Value = "1.643.97"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5146"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06392"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07783"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.309"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "1.647.99"
$ws.Range("E13").Value = "  -3.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5503"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "0.0₅7790"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").Value = "26.061.19"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "199.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.489"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.135"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.904"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1222"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.907"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04894"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.314"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.256"
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.386"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9206"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.602"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5623"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "1.117.26"
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.544"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.579"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8125"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₈122"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.782.47"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4546"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05246"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09603"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
